# Insert a new data row at row 13 (shifts existing rows 13:78 down to 14:79,
# matching the author's "weekly" addition of a new price record ahead of the
# existing ones). After the insert, populate the new row 13 with its values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 13..78 down to 14..79, leaving row 13 blank (ready for new data).
$ws.Rows.Item(13).Insert()

# Populate the newly inserted row 13 with the new record's data.
$ws.Cells.Item(13, 1).Value = 11
$ws.Cells.Item(13, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(13, 3).Value = "Bíobío"
$ws.Cells.Item(13, 4).Value = 44623
$ws.Cells.Item(13, 5).Value = 8
$ws.Cells.Item(13, 6).Value = 100112021
$ws.Cells.Item(13, 7).Value = "Ají"
$ws.Cells.Item(13, 8).Value = "Inferno"
$ws.Cells.Item(13, 9).Value = "Primera"
$ws.Cells.Item(13, 10).Value = 220
$ws.Cells.Item(13, 11).Value = 12000
$ws.Cells.Item(13, 12).Value = 15000
$ws.Cells.Item(13, 13).Value = 13636
$ws.Cells.Item(13, 14).Value = "`$/caja 15 kilos"
$ws.Cells.Item(13, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(13, 16).Value = 909
$ws.Cells.Item(13, 17).Value = 15
$ws.Cells.Item(13, 18).Value = "Hortaliza"

# Ensure the date cell keeps a date number format (matches style used by the
# other rows' Fecha column), and is stored as a true date serial, not text.
$ws.Cells.Item(13, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
